$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$data = @{
    23 = @(14, 9, 17, 5, 16, 3)
    24 = @(5, 1, 13, 4, 3, 1)
    25 = @(14, 14, 14, 9, 13, 13)
    26 = @(13, 3, 10, 10, 10, 0)
    27 = @(8, 0, 4, 8, 14, 14)
    28 = @(6, 7, 7, 14, 5, 4)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 3 + $i   # Column C = 3
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$ws.Range("H28").Select()
